# fix bug add and upload users
# Insert a new "Username" column between "Email" (C) and "Password" (old D),
# shifting Password/Phone/Category one column to the right (D->E, E->F, F->G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; existing D:F (Password/Phone/Category) and their
# data validations shift right to E:G automatically.
$ws.Columns("D:D").Insert()

# Set the header for the newly inserted column.
$ws.Range("D1").Value = "Username"
